# ECOSCOM-4470 - Change records Iterable strategy for group actions
#
# The report template's sample "value" placeholder (A2) is no longer a
# fixed shared string - it's cleared out so the cell keeps its existing
# style/formatting but carries no literal content. This also drops the
# now-unused "value" entry from the shared-strings table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents()
